# Applies the "Knärot" BILAGA 1 section addition + header date update
# described by the source commit diff.
#
# NOTE: this interpreter only reliably binds *positional* parameters to
# user-defined functions (named "-Param value" binding does not stick),
# so all helper functions below take positional args.

$d = $word.ActiveDocument

function New-BodyParagraph($doc, $StyleName) {
    # Inserts a brand-new paragraph right after the current last paragraph
    # of the main document body and applies the requested style to it.
    $count = $doc.Paragraphs.Count
    $last = $doc.Paragraphs.Item($count)
    $lastRange = $last.Range
    $lastRange.Collapse(0)   # wdCollapseEnd
    $lastRange.InsertParagraphAfter()
    $newCount = $doc.Paragraphs.Count
    $newPara = $doc.Paragraphs.Item($newCount)
    $newPara.Style = $StyleName
    return $newPara
}

function Add-Runs($doc, $Para, $Runs) {
    # Writes $Runs (an ordered list of @{ Text = ...; Italic = $true/$false })
    # into $Para as consecutive Word runs, toggling italics per run so that
    # the resulting OOXML keeps each run separate (matching the source diff).
    $fullText = [string]::Join("", ($Runs | ForEach-Object { $_.Text }))
    $r = $Para.Range
    $r.Collapse(0)
    $r.InsertAfter($fullText)
    $paraStart = $Para.Range.Start
    $offset = 0
    foreach ($run in $Runs) {
        $len = $run.Text.Length
        if ($run.Italic -and $len -gt 0) {
            $subStart = $paraStart + $offset
            $subEnd = $paraStart + $offset + $len
            $sub = $doc.Range($subStart, $subEnd)
            $sub.Font.Italic = $true
        }
        $offset += $len
    }
}

# New paragraph 0: style=Heading1
$para0 = New-BodyParagraph $d 'Heading1'
$r0 = $para0.Range
$r0.Collapse(0)
$r0.InsertAfter('Knärot – ekologi samt krav på livsmiljön')

# New paragraph 1: style=Normal
$para1 = New-BodyParagraph $d 'Normal'
$r1 = $para1.Range
$r1.Collapse(0)
$r1.InsertAfter('Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).')

# New paragraph 2: style=Normal
$para2 = New-BodyParagraph $d 'Normal'
Add-Runs $d $para2 @(@{ Text = 'Samuel Johnsons doktorsavhandling '; Italic = $false }, @{ Text = '“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“'; Italic = $true }, @{ Text = ' (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: '; Italic = $false }, @{ Text = '“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” '; Italic = $true }, @{ Text = 'Vidare '; Italic = $false }, @{ Text = '“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”'; Italic = $true })

# New paragraph 3: style=Normal
$para3 = New-BodyParagraph $d 'Normal'
Add-Runs $d $para3 @(@{ Text = 'Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: '; Italic = $false }, @{ Text = '“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”'; Italic = $true })

# New paragraph 4: style=Normal
$para4 = New-BodyParagraph $d 'Normal'
$r4 = $para4.Range
$r4.Collapse(0)
$r4.InsertAfter('En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).')

# New paragraph 5: style=Normal
$para5 = New-BodyParagraph $d 'Normal'
$r5 = $para5.Range
$r5.Collapse(0)
$r5.InsertAfter('Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).')

# New paragraph 6: style=Heading2
$para6 = New-BodyParagraph $d 'Heading2'
$r6 = $para6.Range
$r6.Collapse(0)
$r6.InsertAfter('Referenser - knärot')

# New paragraph 7: style=Normal
$para7 = New-BodyParagraph $d 'Normal'
Add-Runs $d $para7 @(@{ Text = 'de Graaf M & Roberts M.R., 2009. '; Italic = $false }, @{ Text = 'Short-term response of the herbaceous layer within leave patches after harvest. '; Italic = $true }, @{ Text = 'Forest Ecology and Management 257, 1014-1025'; Italic = $false })

# New paragraph 8: style=Normal
$para8 = New-BodyParagraph $d 'Normal'
Add-Runs $d $para8 @(@{ Text = 'Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. '; Italic = $false }, @{ Text = 'Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. '; Italic = $true }, @{ Text = 'Ecological Applications, 22, 2049-2064 '; Italic = $false })

# New paragraph 9: style=Normal
$para9 = New-BodyParagraph $d 'Normal'
Add-Runs $d $para9 @(@{ Text = 'Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. '; Italic = $false }, @{ Text = 'Interactive effects of drought and edge exposure on old-growth forest understory species. '; Italic = $true }, @{ Text = 'Landscape Ecology, 37, sid 1839-1853'; Italic = $false })

# New paragraph 10: style=Normal
$para10 = New-BodyParagraph $d 'Normal'
Add-Runs $d $para10 @(@{ Text = 'Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. '; Italic = $false }, @{ Text = 'Biological legacies buffer local species extinction after logging. '; Italic = $true }, @{ Text = 'Journal of Applied Ecology. 51, 53-62.'; Italic = $false })

# New paragraph 11: style=Normal
$para11 = New-BodyParagraph $d 'Normal'
Add-Runs $d $para11 @(@{ Text = 'Skogsstyrelsen, 2022. '; Italic = $false }, @{ Text = 'Vägledning för hänsyn till knärot. '; Italic = $true }, @{ Text = 'https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/'; Italic = $false })

# New paragraph 12: style=Normal
$para12 = New-BodyParagraph $d 'Normal'
Add-Runs $d $para12 @(@{ Text = 'SLU Artdatabanken, 2021. '; Italic = $false }, @{ Text = 'Artfaktablad. Naturvård – artfakta. '; Italic = $true }, @{ Text = 'SLU Artdatabanken, Uppsala '; Italic = $false })

# Update the "first page" header date from 2023-09-13 to 2023-09-15.
# Headers/footers are separate stories, not reachable via $d.Content, so
# walk every section's headers/footers and replace wherever it is found.
for ($si = 1; $si -le $d.Sections.Count; $si++) {
    $sec = $d.Sections.Item($si)
    for ($hi = 1; $hi -le $sec.Headers.Count; $hi++) {
        $h = $sec.Headers.Item($hi)
        if ($h.Exists) {
            $h.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2) | Out-Null
        }
    }
    for ($fi = 1; $fi -le $sec.Footers.Count; $fi++) {
        $ft = $sec.Footers.Item($fi)
        if ($ft.Exists) {
            $ft.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2) | Out-Null
        }
    }
}
